$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (G=5505)
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

# Row 51 (G=5486)
$ws.Range("H51").Value = 8106.9287
$ws.Range("I51").Value = 2499.5
$ws.Range("J51").Value = 9041.5
$ws.Range("K51").Value = 2499.5
$ws.Range("L51").Value = 9041.5
$ws.Range("M51").Value = -2015.5
$ws.Range("N51").Value = -10009.5

# Row 70 (G=12604)
$ws.Range("H70").Value = 6442.6
$ws.Range("I70").Value = 6249.75
$ws.Range("J70").Value = 6571.1665
$ws.Range("K70").Value = 18749.25
$ws.Range("L70").Value = 19713.4995
$ws.Range("M70").Value = -18479.25
$ws.Range("N70").Value = -20253.4995

# Row 73 (G=12604)
$ws.Range("H73").Value = 6442.6
$ws.Range("I73").Value = 6249.75
$ws.Range("J73").Value = 6571.1665
$ws.Range("K73").Value = 18749.25
$ws.Range("L73").Value = 19713.4995
$ws.Range("M73").Value = -17813.25
$ws.Range("N73").Value = -21585.4995

# Row 96 (G=19894)
$ws.Range("H96").Value = 3268.75
$ws.Range("I96").Value = 1830
$ws.Range("J96").Value = 5666.6665
$ws.Range("K96").Value = 5490
$ws.Range("L96").Value = 16999.9995
$ws.Range("M96").Value = -4117
$ws.Range("N96").Value = -19745.9995

# Row 125 (G=36228)
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -15540

# Row 138 (G=44169)
$ws.Range("H138").Value = 2499.5
$ws.Range("I138").Value = 1000
$ws.Range("J138").Value = 2999.3333
$ws.Range("K138").Value = 3000
$ws.Range("L138").Value = 8997.999899999999
$ws.Range("M138").Value = 2140
$ws.Range("N138").Value = -19277.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (G=43999)
$ws.Range("H61").Value = 10364
$ws.Range("I61").Value = 10364
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10364
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -10152
$ws.Range("N61").ClearContents()

# Row 74 (G=44000)
$ws.Range("H74").Value = 5279.353
$ws.Range("I74").Value = 5333.1665
$ws.Range("J74").Value = 5250
$ws.Range("K74").Value = 5333.1665
$ws.Range("L74").Value = 5250
$ws.Range("M74").Value = -4459.1665
$ws.Range("N74").Value = -6998

# Row 77 (G=44000)
$ws.Range("H77").Value = 5279.353
$ws.Range("I77").Value = 5333.1665
$ws.Range("J77").Value = 5250
$ws.Range("K77").Value = 26665.8325
$ws.Range("L77").Value = 26250
$ws.Range("M77").Value = -22297.8325
$ws.Range("N77").Value = -34986

# Row 136 (G=43999)
$ws.Range("H136").Value = 10364
$ws.Range("I136").Value = 10364
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 31092
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -28542
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 8 (G=2507)
$ws.Range("H8").Value = 1005
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1005
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1005
$ws.Range("N8").Value = -1285
$ws.Range("M8").ClearContents()

# Row 11 (G=2481)
$ws.Range("H11").Value = 992
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 992
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 992
$ws.Range("N11").Value = -1272

$ws = $wb.Worksheets.Item("CRP")
# Row 122 (G=36196)
$ws.Range("H122").Value = 948.2222
$ws.Range("I122").Value = 948.2222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2844.6666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -394.6666

$ws = $wb.Worksheets.Item("CUL")
# Row 8 (G=16734)
$ws.Range("H8").Value = 942.6
$ws.Range("I8").Value = 942.6
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2827.8
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2688.8

# Row 68 (G=12895)
$ws.Range("H68").Value = 1499.6666
$ws.Range("I68").Value = 1499.6666
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 4498.9998
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -3687.9998

# Row 71 (G=12895)
$ws.Range("H71").Value = 1499.6666
$ws.Range("I71").Value = 1499.6666
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 13496.9994
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -9440.999400000001

# Row 80 (G=12890)
$ws.Range("H80").Value = 7001.5
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 10003
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 30009
$ws.Range("M80").Value = -11064
$ws.Range("N80").Value = -31881

# Row 83 (G=12890)
$ws.Range("H83").Value = 7001.5
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 10003
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 90027
$ws.Range("M83").Value = -31320
$ws.Range("N83").Value = -99387

$ws = $wb.Worksheets.Item("GSM")
# Row 62 (G=11983)
$ws.Range("H62").Value = 25000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 25000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26372

# Row 65 (G=11983)
$ws.Range("H65").Value = 25000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 25000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81864

# Row 80 (G=12521)
$ws.Range("H80").Value = 12998
$ws.Range("I80").Value = 12998
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 12998
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -12000
$ws.Range("N80").ClearContents()

# Row 83 (G=12521)
$ws.Range("H83").Value = 12998
$ws.Range("I83").Value = 12998
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 64990
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -59998
$ws.Range("N83").ClearContents()

# Row 93 (G=18107)
$ws.Range("H93").Value = 35000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 35000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744

# Row 111 (G=25853)
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 113 (G=27710)
$ws.Range("H113").Value = 2711.3333
$ws.Range("I113").Value = 3106
$ws.Range("J113").Value = 738
$ws.Range("K113").Value = 3106
$ws.Range("L113").Value = 738
$ws.Range("M113").Value = -936
$ws.Range("N113").Value = -5078

# Row 124 (G=34247)
$ws.Range("H124").Value = 27500
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 27500
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 27500
$ws.Range("N124").Value = -37320

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G=5277)
$ws.Range("H22").Value = 4351.448
$ws.Range("I22").Value = 4704.857
$ws.Range("J22").Value = 3423.75
$ws.Range("K22").Value = 4704.857
$ws.Range("L22").Value = 3423.75
$ws.Range("M22").Value = -4409.857
$ws.Range("N22").Value = -4013.75

# Row 27 (G=5277)
$ws.Range("H27").Value = 4351.448
$ws.Range("I27").Value = 4704.857
$ws.Range("J27").Value = 3423.75
$ws.Range("K27").Value = 4704.857
$ws.Range("L27").Value = 3423.75
$ws.Range("M27").Value = -4597.857
$ws.Range("N27").Value = -3637.75

# Row 63 (G=12006)
$ws.Range("H63").Value = 25000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 25000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 25000
$ws.Range("N63").Value = -26498

# Row 66 (G=12006)
$ws.Range("H66").Value = 25000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 25000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 75000
$ws.Range("N66").Value = -82488

# Row 132 (G=44058)
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Row 133 (G=41903)
$ws.Range("H133").Value = 30000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 30000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -35060

# Row 136 (G=44060)
$ws.Range("H136").Value = 14998
$ws.Range("I136").Value = 14998
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 44994
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -42444
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (G=27746)
$ws.Range("H107").Value = 1194.5
$ws.Range("I107").Value = 1194.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3583.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1663.5

# Row 136 (G=44031)
$ws.Range("H136").Value = 7131.7036
$ws.Range("I136").Value = 6502.95
$ws.Range("J136").Value = 8928.143
$ws.Range("K136").Value = 19508.85
$ws.Range("L136").Value = 26784.429
$ws.Range("M136").Value = -16958.85
$ws.Range("N136").Value = -31884.429
